# Add data for 2022-06-05
# - Rename sheet / "through" labels from 2022-05-27 to 2022-05-28
# - Update the May row (row 6) and Total row (row 7) with the new daily counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Through 2022-05-28"

# Update the "May (through ...)" label in column A, row 6
$ws.Range("A6").Value = "May (through 05-28)"

# Update the May row (row 6) values for years 2015-2022 (columns B-I)
$ws.Range("B6").Value = 17
$ws.Range("C6").Value = 42
$ws.Range("D6").Value = 55
$ws.Range("E6").Value = 46
$ws.Range("F6").Value = 42
$ws.Range("G6").Value = 61
$ws.Range("H6").Value = 102
$ws.Range("I6").Value = 101

# Update the Total row (row 7) values for years 2015-2022 (columns B-I)
$ws.Range("B7").Value = 106
$ws.Range("C7").Value = 204
$ws.Range("D7").Value = 308
$ws.Range("E7").Value = 292
$ws.Range("F7").Value = 197
$ws.Range("G7").Value = 323
$ws.Range("H7").Value = 625
$ws.Range("I7").Value = 652
